$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:41:15.382154",
    "2021-10-05 13:41:15.382167",
    "2021-10-05 13:41:15.382171",
    "2021-10-05 13:41:15.382174",
    "2021-10-05 13:41:15.382178",
    "2021-10-05 13:41:15.382181",
    "2021-10-05 13:41:15.382184",
    "2021-10-05 13:41:15.382187",
    "2021-10-05 13:41:15.382190",
    "2021-10-05 13:41:15.382193",
    "2021-10-05 13:41:15.382196",
    "2021-10-05 13:41:15.382199",
    "2021-10-05 13:41:15.382202",
    "2021-10-05 13:41:15.382205",
    "2021-10-05 13:41:15.382208",
    "2021-10-05 13:41:15.382211",
    "2021-10-05 13:41:15.382215",
    "2021-10-05 13:41:15.382218"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
